# Applies the create_detector script tweaks to the intersection setup sheet.
# Updates Detector_Distance (H) and Detect_Actual_Distance (R) values for a
# handful of rows, plus the actual-lane/edge relabeling (N/P) for two rows
# whose source edge changed (Campus_EB -> gneE2, Campus_WB -> gneE1.93).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bryce_NB_0
$ws.Range("H2").Value = 0
$ws.Range("R2").Value = 39.27

# Row 3: Bryce_NB_1
$ws.Range("H3").Value = 0
$ws.Range("R3").Value = 39.27

# Row 5: Campus_EB_0 -> gneE2_0
$ws.Range("H5").Value = 200
$ws.Range("N5").Value = "gneE2_0"
$ws.Range("P5").Value = "gneE2"
$ws.Range("R5").Value = 152.8368118543222

# Row 7: Bryce_SB_0
$ws.Range("H7").Value = 0
$ws.Range("R7").Value = 71.92

# Row 8: Bryce_SB_1
$ws.Range("H8").Value = 0
$ws.Range("R8").Value = 71.92

# Row 11: Campus_WB_1 -> gneE1.93_1
$ws.Range("H11").Value = 200
$ws.Range("N11").Value = "gneE1.93_1"
$ws.Range("P11").Value = "gneE1.93"
$ws.Range("R11").Value = 11.98563409193079
